$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading "aarathy.nair@quadance.com;" address from the
# Recipient Email cell (column P, row 2).
$ws.Range("P2").Value = "yedu.yesodharan@quadance.com;mafna.janeefar@quadance.com"

# Update the view: scroll the pane so column O is left-most visible,
# and move the active selection to P10.
$win = $excel.ActiveWindow
$win.ScrollColumn = 15
$win.ScrollRow = 1
$ws.Range("P10").Select()
